$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25, columns B through N (column A = time step index is unchanged,
# column G stays 1, column H has no data in this sheet).
$data = @(
    ,(1.02,1.034611238139842,1.049632132427115,1.03353711852611,1.056232926522912,1.041589144272035,1.039729534889943,1.052388650759266,1.036339182953833,1.058971200727039,1.041206070088131)
    ,(1.02,1.037385856187176,1.050594487172144,1.035960197764183,1.057705768819247,1.041917424593009,1.042138146417881,1.053162903096879,1.038566961628376,1.060255963414404,1.043618102120711)
    ,(1.02,1.039170298156628,1.051212606778709,1.037518240930108,1.058651968009541,1.042125867074013,1.043686020991715,1.053658699461481,1.039998372205912,1.061079966342351,1.045168174853983)
    ,(1.02,1.039917933717811,1.051471380249789,1.038170944648037,1.059048142669023,1.042212551523705,1.044334257741408,1.053865902210756,1.040597774245022,1.061424649704589,1.045817332174215)
    ,(1.02,1.040043317881144,1.05151476630122,1.038280403487227,1.059114568656077,1.042227051031581,1.044442955324237,1.05390062073692,1.040698279577637,1.061482423005363,1.045926184120078)
    ,(1.02,1.039180298001321,1.05121606876255,1.037526971336654,1.05865726799468,1.042127029058372,1.043694692463737,1.053661472929154,1.040006390653637,1.061084578771354,1.04517685864049)
    ,(1.02,1.035551237685792,1.049958323541126,1.034358086175638,1.056732110876182,1.041700916021763,1.040545780885899,1.052651397528439,1.03709419984215,1.05940692470392,1.042023475247026)
    ,(1.02,1.029069491965422,1.047706212346415,1.028695883184464,1.053286195082025,1.040919256915294,1.034912541661665,1.050831069944601,1.031882473489659,1.056393395356999,1.036382236176952)
    ,(1.02,1.024685320151746,1.046179783280677,1.024864544014617,1.050951139628743,1.040376965241682,1.031096247068501,1.049589400662865,1.028350405815758,1.054344156473794,1.032560522007731)
    ,(1.02,1.022770934997472,1.045512669818372,1.023191224607576,1.049930667763367,1.040137011075652,1.029428413336057,1.049044862891133,1.026806471111656,1.053446884807793,1.030890319760083)
    ,(1.02,1.022057350309824,1.045263929496042,1.022567448686394,1.04955017446352,1.040047099600238,1.028806518939806,1.048841543908645,1.026230728024197,1.05311207088931,1.030267542202458)
    ,(1.02,1.022210531130234,1.04531732820411,1.022701352962323,1.049631857525134,1.040066421457113,1.028940026711483,1.048885204496806,1.026354330225155,1.053183959284314,1.030401239570484)
    ,(1.02,1.022712001336975,1.045492128250458,1.023139709060631,1.049899245748671,1.040129594971866,1.02937705651007,1.049028078080314,1.026758926489387,1.053419240342333,1.030838890001513)
    ,(1.02,1.023020640052837,1.045599702563054,1.023409496761555,1.05006379994662,1.040168414361838,1.029646006302177,1.049115967053606,1.027007910262099,1.053564001348772,1.031108221733174)
    ,(1.02,1.024812026524752,1.04622392599634,1.024975288260458,1.051018664383699,1.040392781117946,1.031206605443556,1.049625393156946,1.028452559304804,1.054403493161776,1.032671037104337)
    ,(1.02,1.025931365794907,1.046613822335002,1.025953578088348,1.051615090071054,1.040532137931939,1.032181361750759,1.049943085840382,1.029354808370051,1.054927399711763,1.033647177677302)
    ,(1.02,1.02658271903533,1.046840649290535,1.026522821314141,1.051962072760561,1.040612927313162,1.032748444783746,1.050127727427439,1.029879678562693,1.055232029054373,1.034215066033241)
    ,(1.02,1.026804555453368,1.046917891415959,1.02671668738469,1.052080232910168,1.040640390722227,1.03294155751277,1.0501905735101,1.030058411109613,1.055335738575124,1.034408453004491)
    ,(1.02,1.025811431064613,1.046572051620493,1.02584875980622,1.051551192826538,1.040517237544548,1.032076932848499,1.049909069151916,1.029258150556176,1.054871288675624,1.033542600474023)
    ,(1.02,1.022564400530429,1.0454406802288,1.023010686331178,1.049820546751582,1.040111013608704,1.029248428664589,1.048986034603416,1.026639845886631,1.053349998374689,1.030710079489737)
    ,(1.02,1.020508370761497,1.044723866561316,1.021213331017096,1.048724044376404,1.03985107607353,1.027456185998229,1.048399581773544,1.024980517177038,1.052384648822325,1.028915291633219)
    ,(1.02,1.021599715975566,1.045104388648751,1.022167396952022,1.049306127221053,1.039989306484678,1.028407627835758,1.048711056539813,1.025861425555287,1.052897249817965,1.029868084627289)
    ,(1.02,1.025865629155432,1.046590927839947,1.025896126925768,1.051580068027459,1.040523971917644,1.032124124330322,1.049924441877635,1.029301830372585,1.054896645791516,1.033589858973163)
    ,(1.02,1.03075594031525,1.048292776036759,1.030169368341855,1.05418357229829,1.041125032434374,1.036379285032658,1.051306554804023,1.033239698882017,1.057179426317301,1.037851062491769)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 6).Value = $vals[4]   # F
    # column G (7) stays 1, unchanged
    # column H (8) has no data, left untouched
    $ws.Cells.Item($row, 9).Value = $vals[5]   # I
    $ws.Cells.Item($row, 10).Value = $vals[6]  # J
    $ws.Cells.Item($row, 11).Value = $vals[7]  # K
    $ws.Cells.Item($row, 12).Value = $vals[8]  # L
    $ws.Cells.Item($row, 13).Value = $vals[9]  # M
    $ws.Cells.Item($row, 14).Value = $vals[10] # N
}
